$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 899.6667
$ws.Range("I4").Value = 899.6667
$ws.Range("K4").Value = 899.6667
$ws.Range("M4").Value = -785.6667
$ws.Range("H28").Value = 349.86667
$ws.Range("I28").Value = 140
$ws.Range("J28").Value = 927
$ws.Range("K28").Value = 140
$ws.Range("L28").Value = 927
$ws.Range("M28").Value = 345
$ws.Range("N28").Value = -1897
$ws.Range("H62").Value = 18535.143
$ws.Range("I62").Value = 19980.637
$ws.Range("K62").Value = 19980.637
$ws.Range("M62").Value = -19356.637
$ws.Range("H65").Value = 18535.143
$ws.Range("I65").Value = 19980.637
$ws.Range("K65").Value = 99903.185
$ws.Range("M65").Value = -96783.185
$ws.Range("H76").Value = 5765.375
$ws.Range("I76").Value = 3964
$ws.Range("J76").Value = 7166.4443
$ws.Range("K76").Value = 3964
$ws.Range("L76").Value = 7166.4443
$ws.Range("M76").Value = -3649
$ws.Range("N76").Value = -7796.4443
$ws.Range("H79").Value = 5765.375
$ws.Range("I79").Value = 3964
$ws.Range("J79").Value = 7166.4443
$ws.Range("K79").Value = 3964
$ws.Range("L79").Value = 7166.4443
$ws.Range("M79").Value = -2872
$ws.Range("N79").Value = -9350.444299999999
$ws.Range("H109").Value = 100570
$ws.Range("J109").Value = 100570
$ws.Range("L109").Value = 100570
$ws.Range("N109").Value = -103344
$ws.Range("H111").Value = 4492.8887
$ws.Range("I111").Value = 2587.6
$ws.Range("K111").Value = 7762.799999999999
$ws.Range("M111").Value = -4695.799999999999
$ws.Range("H137").Value = 1162095.4
$ws.Range("I137").Value = 1531
$ws.Range("J137").Value = 2322659.8
$ws.Range("K137").Value = 4593
$ws.Range("L137").Value = 6967979.399999999
$ws.Range("M137").Value = -2043
$ws.Range("N137").Value = -6973079.399999999
$ws.Range("H138").Value = 4413.3486
$ws.Range("I138").Value = 1769.7142
$ws.Range("J138").Value = 5689.5864
$ws.Range("K138").Value = 5309.142599999999
$ws.Range("L138").Value = 17068.7592
$ws.Range("M138").Value = -169.1425999999992
$ws.Range("N138").Value = -27348.7592
$ws.Range("H141").Value = 2246.1667
$ws.Range("I141").Value = 2187.7273
$ws.Range("K141").Value = 6563.1819
$ws.Range("M141").Value = -1383.1819

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 13331.556
$ws.Range("I2").Value = 5718.091
$ws.Range("K2").Value = 5718.091
$ws.Range("M2").Value = -5605.091
$ws.Range("H32").Value = 2675.1594
$ws.Range("I32").Value = 2749.111
$ws.Range("K32").Value = 2749.111
$ws.Range("M32").Value = -2462.111
$ws.Range("H57").Value = 11122774
$ws.Range("I57").Value = 11122774
$ws.Range("K57").Value = 11122774
$ws.Range("M57").Value = -11122290
$ws.Range("H116").Value = 13331.556
$ws.Range("I116").Value = 5718.091
$ws.Range("K116").Value = 5718.091
$ws.Range("M116").Value = -3424.091
$ws.Range("H122").Value = 3595.6086
$ws.Range("I122").Value = 3164.2727
$ws.Range("J122").Value = 3991
$ws.Range("K122").Value = 9492.8181
$ws.Range("L122").Value = 11973
$ws.Range("M122").Value = -7042.8181
$ws.Range("N122").Value = -16873

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 13331.556
$ws.Range("I3").Value = 5718.091
$ws.Range("K3").Value = 5718.091
$ws.Range("M3").Value = -5604.091
$ws.Range("H94").Value = 1139.6364
$ws.Range("I94").Value = 1139.6364
$ws.Range("K94").Value = 1139.6364
$ws.Range("M94").Value = -688.6364000000001
$ws.Range("H107").Value = 2206.8572
$ws.Range("I107").Value = 2346.9443
$ws.Range("J107").Value = 1366.3334
$ws.Range("K107").Value = 2346.9443
$ws.Range("L107").Value = 1366.3334
$ws.Range("M107").Value = -426.9443000000001
$ws.Range("N107").Value = -5206.3334
$ws.Range("H134").Value = 2406491.8
$ws.Range("I134").Value = 2939888
$ws.Range("K134").Value = 8819664
$ws.Range("M134").Value = -8817129

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 1671.7778
$ws.Range("I12").Value = 839.1667
$ws.Range("J12").Value = 3337
$ws.Range("K12").Value = 839.1667
$ws.Range("L12").Value = 3337
$ws.Range("M12").Value = -669.1667
$ws.Range("N12").Value = -3677
$ws.Range("H22").Value = 997.64703
$ws.Range("I22").Value = 758.4545000000001
$ws.Range("J22").Value = 1436.1666
$ws.Range("K22").Value = 758.4545000000001
$ws.Range("L22").Value = 1436.1666
$ws.Range("M22").Value = -408.4545000000001
$ws.Range("N22").Value = -2136.1666
$ws.Range("H31").Value = 298397.4
$ws.Range("I31").Value = 518147.16
$ws.Range("K31").Value = 518147.16
$ws.Range("M31").Value = -517852.16
$ws.Range("H34").Value = 298397.4
$ws.Range("I34").Value = 518147.16
$ws.Range("K34").Value = 518147.16
$ws.Range("M34").Value = -517945.16
$ws.Range("H134").Value = 5441.7856
$ws.Range("I134").Value = 6158.8335
$ws.Range("K134").Value = 18476.5005
$ws.Range("M134").Value = -15941.5005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 5030.1777
$ws.Range("I68").Value = 1575
$ws.Range("J68").Value = 5367.268
$ws.Range("K68").Value = 4725
$ws.Range("L68").Value = 16101.804
$ws.Range("M68").Value = -3914
$ws.Range("N68").Value = -17723.804
$ws.Range("H71").Value = 5030.1777
$ws.Range("I71").Value = 1575
$ws.Range("J71").Value = 5367.268
$ws.Range("K71").Value = 14175
$ws.Range("L71").Value = 48305.412
$ws.Range("M71").Value = -10119
$ws.Range("N71").Value = -56417.412
$ws.Range("H107").Value = 889.5333000000001
$ws.Range("I107").Value = 842.5
$ws.Range("J107").Value = 983.6
$ws.Range("K107").Value = 2527.5
$ws.Range("L107").Value = 2950.8
$ws.Range("M107").Value = -607.5
$ws.Range("N107").Value = -6790.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 44500
$ws.Range("I5").Value = 44000
$ws.Range("J5").Value = 45000
$ws.Range("K5").Value = 44000
$ws.Range("L5").Value = 45000
$ws.Range("M5").Value = -43888
$ws.Range("N5").Value = -45224
$ws.Range("H107").Value = 19719.111
$ws.Range("I107").Value = 25231.572
$ws.Range("J107").Value = 425.5
$ws.Range("K107").Value = 25231.572
$ws.Range("L107").Value = 425.5
$ws.Range("M107").Value = -23311.572
$ws.Range("N107").Value = -4265.5
$ws.Range("H113").Value = 5367.951
$ws.Range("I113").Value = 4816.6665
$ws.Range("J113").Value = 6146.2354
$ws.Range("K113").Value = 4816.6665
$ws.Range("L113").Value = 6146.2354
$ws.Range("M113").Value = -2646.6665
$ws.Range("N113").Value = -10486.2354

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H39").Value = 2000
$ws.Range("I39").Value = 2000
$ws.Range("K39").Value = 2000
$ws.Range("M39").Value = -1540
$ws.Range("H61").Value = 3251.0908
$ws.Range("I61").Value = 3045.375
$ws.Range("J61").Value = 3799.6667
$ws.Range("K61").Value = 3045.375
$ws.Range("L61").Value = 3799.6667
$ws.Range("M61").Value = -2843.375
$ws.Range("N61").Value = -4203.6667
$ws.Range("H68").Value = 2759.1428
$ws.Range("I68").Value = 2675.4443
$ws.Range("J68").Value = 2909.8
$ws.Range("K68").Value = 2675.4443
$ws.Range("L68").Value = 2909.8
$ws.Range("M68").Value = -1926.4443
$ws.Range("N68").Value = -4407.8
$ws.Range("H71").Value = 2759.1428
$ws.Range("I71").Value = 2675.4443
$ws.Range("J71").Value = 2909.8
$ws.Range("K71").Value = 13377.2215
$ws.Range("L71").Value = 14549
$ws.Range("M71").Value = -9633.2215
$ws.Range("N71").Value = -22037
$ws.Range("H113").Value = 3251.0908
$ws.Range("I113").Value = 3045.375
$ws.Range("J113").Value = 3799.6667
$ws.Range("K113").Value = 3045.375
$ws.Range("L113").Value = 3799.6667
$ws.Range("M113").Value = -875.375
$ws.Range("N113").Value = -8139.6667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1293.7273
$ws.Range("I81").Value = 1205
$ws.Range("J81").Value = 1400.2
$ws.Range("K81").Value = 2410
$ws.Range("L81").Value = 2800.4
$ws.Range("M81").Value = -1349
$ws.Range("N81").Value = -4922.4
$ws.Range("H84").Value = 1293.7273
$ws.Range("I84").Value = 1205
$ws.Range("J84").Value = 1400.2
$ws.Range("K84").Value = 12050
$ws.Range("L84").Value = 14002
$ws.Range("M84").Value = -6746
$ws.Range("N84").Value = -24610

